$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# --- Split TC003 ("Savings or Chequing") into two rows: TC003 (Savings) and a new TC004 (Chequing) ---
# Insert a new row at position 5; this shifts the old rows 5-8 (TC004..TC007) down to 6-9.
$ws.Rows.Item(5).Insert()

# Row 4 keeps being TC003 but now only covers the Savings case.
$ws.Cells.Item(4, 2).Value2  = "TC003"
$ws.Cells.Item(4, 3).Value2  = "Verify that an account can be created for Savings"
$ws.Cells.Item(4, 5).Value2  = "1. User knows what kind of account they want to create"
$ws.Cells.Item(4, 6).Value2  = "1. Create a Savings Account using the createAcc() method of CreateAccount class"
$ws.Cells.Item(4, 7).Value2  = "1. A New Savings Account object is created.                   "
$ws.Cells.Item(4, 8).Value2  = "1. If a Savings Account is created, a message saying a Savings Account has been created should display."
$ws.Cells.Item(4, 9).Value2  = "1. If Savings is inserted for card type, the result is as specified.                           "
$ws.Rows.Item(4).RowHeight = 60

# New row 5 becomes the Chequing counterpart, re-using TC004 as its id (later ids shift by one).
$ws.Cells.Item(5, 1).Value2  = "TS01"
$ws.Cells.Item(5, 2).Value2  = "TC004"
$ws.Cells.Item(5, 3).Value2  = "Verify that an account can be created for Chequing"
$ws.Cells.Item(5, 4).Value2  = "N/A"
$ws.Cells.Item(5, 5).Value2  = "1. User knows what kind of account they want to create"
$ws.Cells.Item(5, 6).Value2  = "                                                        1. Create a Chequing Account using the createAcc() method of CreateAccount class"
$ws.Cells.Item(5, 7).Value2  = "1. A New Chequings Account object Is created."
$ws.Cells.Item(5, 8).Value2  = "1. If a Chequings Account is created, a message saying a Chequings Account has been created should display."
$ws.Cells.Item(5, 9).Value2  = "1. If Chequing is inserted for card type, the result is as specified."
$ws.Cells.Item(5, 10).Value2 = "Pass"
$ws.Cells.Item(5, 11).Value2 = "Create Accounts Test Case"
$ws.Cells.Item(5, 12).Value2 = "Tyler Serio"
$ws.Cells.Item(5, 13).Value2 = "03/24/2015"
$ws.Cells.Item(5, 14).Value2 = "Tyler Serio"
$ws.Cells.Item(5, 15).Value2 = "04/13/2015"
$ws.Cells.Item(5, 16).Value2 = "OS: Windows 8.1                   IDE: Eclipse"
$ws.Rows.Item(5).RowHeight = 60

# --- Renumber the ids of the rows that got pushed down by the insert (old TC004..TC007 -> TC005..TC008) ---
$ws.Cells.Item(6, 2).Value2 = "TC005"
$ws.Cells.Item(7, 2).Value2 = "TC006"
$ws.Cells.Item(8, 2).Value2 = "TC007"
$ws.Cells.Item(9, 2).Value2 = "TC008"

# --- Update the sheet's view to match the new selection/scroll position ---
$ws.Range("A9").Select()
$window = $excel.ActiveWindow
$window.ScrollRow = 9
$window.ScrollColumn = 1
$ws.Range("B9").Select()
